# Update the LR-pairs table (rows 2-21, cols A:T) with refreshed TPM-based
# NATMI ligand-receptor statistics (new sending-cluster mix + recomputed
# detection rates / specificities / weights).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Fgf1", "Nrp1", "ECs", 2, 0.6666666666666666, 0.265349, 0.7960469999999999, 0.1498685997319469, 0.1498685997319469, 3, 1, 86.89540866666668, 260.686226, 0.319779657009892, 0.3197796570098919, 23.05760979429134, 207.518488148622, 0.04792492941883476, 0.04792492941883474),
    @("ECs", "Fgf1", "Nrp1", "FAPs", 2, 0.6666666666666666, 0.265349, 0.7960469999999999, 0.1498685997319469, 0.1498685997319469, 3, 1, 54.34621166666667, 163.038635, 0.1999969065479545, 0.1999969065479545, 14.42071291953833, 129.786416275845, 0.02997325633506298, 0.02997325633506298),
    @("ECs", "Fgf1", "Nrp1", "Inflammatory-Mac", 2, 0.6666666666666666, 0.265349, 0.7960469999999999, 0.1498685997319469, 0.1498685997319469, 3, 1, 60.92601633333334, 182.778049, 0.224210932487692, 0.224210932487692, 16.16665750803367, 145.499917572303, 0.03360217849652447, 0.03360217849652447),
    @("ECs", "Fgf1", "Nrp1", "MuSCs", 2, 0.6666666666666666, 0.265349, 0.7960469999999999, 0.1498685997319469, 0.1498685997319469, 3, 1, 7.809668333333332, 23.429005, 0.02873998867505581, 0.02873998867505581, 2.072287682581667, 18.650589143235, 0.004307221859042625, 0.004307221859042625),
    @("ECs", "Fgf1", "Nrp1", "Resolving-Mac", 2, 0.6666666666666666, 0.265349, 0.7960469999999999, 0.1498685997319469, 0.1498685997319469, 3, 1, 61.75795633333333, 185.273869, 0.2272725152794058, 0.2272725152794058, 16.38741195509367, 147.486707595843, 0.03406101362248205, 0.03406101362248205),
    @("FAPs", "Fgf1", "Nrp1", "ECs", 3, 1, 1.003400666666667, 3.010202, 0.5667187473230926, 0.5667187473230925, 3, 1, 86.89540866666668, 260.686226, 0.319779657009892, 0.3197796570098919, 87.19091098640578, 784.718198877652, 0.1812251266400542, 0.1812251266400541),
    @("FAPs", "Fgf1", "Nrp1", "FAPs", 3, 1, 1.003400666666667, 3.010202, 0.5667187473230926, 0.5667187473230925, 3, 1, 54.34621166666667, 163.038635, 0.1999969065479545, 0.1999969065479545, 54.53102501714111, 490.7792251542699, 0.1133419963473504, 0.1133419963473504),
    @("FAPs", "Fgf1", "Nrp1", "Inflammatory-Mac", 3, 1, 1.003400666666667, 3.010202, 0.5667187473230926, 0.5667187473230925, 3, 1, 60.92601633333334, 182.778049, 0.224210932487692, 0.224210932487692, 61.13320540621089, 550.1988486558979, 0.1270645387955673, 0.1270645387955673),
    @("FAPs", "Fgf1", "Nrp1", "MuSCs", 3, 1, 1.003400666666667, 3.010202, 0.5667187473230926, 0.5667187473230925, 3, 1, 7.809668333333332, 23.429005, 0.02873998867505581, 0.02873998867505581, 7.836226412112221, 70.52603770900998, 0.0162874903800075, 0.0162874903800075),
    @("FAPs", "Fgf1", "Nrp1", "Resolving-Mac", 3, 1, 1.003400666666667, 3.010202, 0.5667187473230926, 0.5667187473230925, 3, 1, 61.75795633333333, 185.273869, 0.2272725152794058, 0.2272725152794058, 61.96797455683755, 557.7117710115378, 0.1287995951601133, 0.1287995951601132),
    @("MuSCs", "Fgf1", "Nrp1", "ECs", 3, 1, 0.4790693333333333, 1.437208, 0.2705774288246194, 0.2705774288246194, 3, 1, 86.89540866666668, 260.686226, 0.319779657009892, 0.3197796570098919, 41.62892549966756, 374.660329497008, 0.08652515738415525, 0.08652515738415524),
    @("MuSCs", "Fgf1", "Nrp1", "FAPs", 3, 1, 0.4790693333333333, 1.437208, 0.2705774288246194, 0.2705774288246194, 3, 1, 54.34621166666667, 163.038635, 0.1999969065479545, 0.1999969065479545, 26.03560339234222, 234.32043053108, 0.05411464874662324, 0.05411464874662323),
    @("MuSCs", "Fgf1", "Nrp1", "Inflammatory-Mac", 3, 1, 0.4790693333333333, 1.437208, 0.2705774288246194, 0.2705774288246194, 3, 1, 60.92601633333334, 182.778049, 0.224210932487692, 0.224210932487692, 29.18778602746578, 262.690074247192, 0.06066641762689003, 0.06066641762689003),
    @("MuSCs", "Fgf1", "Nrp1", "MuSCs", 3, 1, 0.4790693333333333, 1.437208, 0.2705774288246194, 0.2705774288246194, 3, 1, 7.809668333333332, 23.429005, 0.02873998867505581, 0.02873998867505581, 3.741372602004443, 33.67235341803999, 0.007776392240145283, 0.007776392240145283),
    @("MuSCs", "Fgf1", "Nrp1", "Resolving-Mac", 3, 1, 0.4790693333333333, 1.437208, 0.2705774288246194, 0.2705774288246194, 3, 1, 61.75795633333333, 185.273869, 0.2272725152794058, 0.2272725152794058, 29.58634296863911, 266.277086717752, 0.06149481282680567, 0.06149481282680566),
    @("Resolving-Mac", "Fgf1", "Nrp1", "ECs", 1, 0.3333333333333333, 0.02272533333333333, 0.068176, 0.01283522412034115, 0.01283522412034115, 3, 1, 86.89540866666668, 260.686226, 0.319779657009892, 0.3197796570098919, 1.974727127086223, 17.772544143776, 0.004104443566847784, 0.004104443566847783),
    @("Resolving-Mac", "Fgf1", "Nrp1", "FAPs", 1, 0.3333333333333333, 0.02272533333333333, 0.068176, 0.01283522412034115, 0.01283522412034115, 3, 1, 54.34621166666667, 163.038635, 0.1999969065479545, 0.1999969065479545, 1.235035775528889, 11.11532197976, 0.00256700511891792, 0.00256700511891792),
    @("Resolving-Mac", "Fgf1", "Nrp1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.02272533333333333, 0.068176, 0.01283522412034115, 0.01283522412034115, 3, 1, 60.92601633333334, 182.778049, 0.224210932487692, 0.224210932487692, 1.384564029847111, 12.461076268624, 0.002877797568710204, 0.002877797568710204),
    @("Resolving-Mac", "Fgf1", "Nrp1", "MuSCs", 1, 0.3333333333333333, 0.02272533333333333, 0.068176, 0.01283522412034115, 0.01283522412034115, 3, 1, 7.809668333333332, 23.429005, 0.02873998867505581, 0.02873998867505581, 0.1774773160977778, 1.59729584488, 0.0003688841958604077, 0.0003688841958604077),
    @("Resolving-Mac", "Fgf1", "Nrp1", "Resolving-Mac", 1, 0.3333333333333333, 0.02272533333333333, 0.068176, 0.01283522412034115, 0.01283522412034115, 3, 1, 61.75795633333333, 185.273869, 0.2272725152794058, 0.2272725152794058, 1.403470143660444, 12.631231292944, 0.002917093670004832, 0.002917093670004831)
)

$numRows = $rows.Count
$numCols = $rows[0].Count
$data = New-Object 'object[,]' $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r,$c] = $rows[$r][$c]
    }
}

$ws.Range("A2:T21").Value = $data
